# "fix field naming issues"
# Row 2 (B2:Y2) holds the internal field-name keys for this form template.
# They were misnamed with an "Hss" infix (leftover from a copy/paste of the
# "Hx HSS" page) and need to read "Nss" instead (this sheet is "Hx NSS").
# The visible question labels in row 3 are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 25   # column Y
for ($c = 2; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(2, $c)
    $current = $cell.Value2
    if ($current -ne $null -and $current -like "hxHssQ*") {
        $cell.Value = $current -replace "^hxHssQ", "hxNssQ"
    }
}

# Reflect where the editor's cursor ended up after making the fix.
$ws.Range("V12").Select() | Out-Null
